$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Compartments")

# Insert 3 new columns before the existing Volume column (D:F), shifting old D..N to G..Q
$ws.Range("D1:F1").EntireColumn.Insert()

# Rename existing "Type" header (C1) and add headers for the 3 new columns
$ws.Range("C1").Value = "Biological type"
$ws.Range("D1").Value = "Physical type"
$ws.Range("E1").Value = "Geometry"
$ws.Range("F1").Value = "Parent compartment"

# Row 2 (c / cytosol)
$ws.Range("C2").Value = "cellular"
$ws.Range("D2").Value = "fluid"
$ws.Range("E2").Value = "3d"
$ws.Range("F2").Value = "e"

# Row 3 (e / extracellular space)
$ws.Range("C3").Value = "extracellular"
$ws.Range("D3").Value = "fluid"
$ws.Range("E3").Value = "3d"

$null = $ws.Range("A1:H3").AutoFilter()
$af = $ws.AutoFilter
Write-Host "AutoFilter:" $af.Range.Address()
Write-Host "done"
for ($r = 1; $r -le 3; $r++) {
    for ($c = 1; $c -le 17; $c++) {
        $v = $ws.Cells.Item($r, $c).Value()
        Write-Host "R$r C$c = $v"
    }
}
